$p = $ppt.ActivePresentation

# --- Edit 1: Slide 12 - "release" -> "releases", add sentence about feedback window ---
$s1 = $p.Slides.Item(12)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$para1 = $tr1.Paragraphs(1,1)
$whole1 = $para1.Characters(1, 107)
$whole1.Text = "9/9/2019 President Fritz releases a new draft on “collegial” governance, and referendum to be held 9/27-10/1. Feed back limited to 8 days."

# --- Edit 2: Slide 16 - fix "**taken8*" typo to "taken" (bold) ---
$s2 = $p.Slides.Item(16)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$para2 = $tr2.Paragraphs(2,1)
$whole2 = $para2.Characters(1, 260)
$whole2.Text = "The membership of the Executive Committee shall be taken from the College Senate membership, one member from each academic division/school, provost, chief of staff, one member from the Higher Education Officer series, and the President who serves as chair."
$bold2 = $para2.Characters(52, 5)
$bold2.Font.Bold = $true

# --- Edit 3: Slide 17 - punctuation updates ---
$s3 = $p.Slides.Item(17)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$para3 = $tr3.Paragraphs(1,1)
$whole3 = $para3.Characters(1, 248)
$whole3.Text = "Note the original selection process … the word used is “taken,” not elected. This was modified, relaxing the faculty membership to be more democratic. but pay attention – it doesn’t say who actually votes! (There is a “from” and a “by” typically.)"

# --- Edit 4: Slide 23 - expand "negative proposal..." sentence ---
$s4 = $p.Slides.Item(23)
$tr4 = $s4.Shapes.Item(1).TextFrame.TextRange
$para4 = $tr4.Paragraphs(3,1)
$whole4 = $para4.Characters(1, 70)
$whole4.Text = "The latter allows a proposal submitted by the president with a negative recommendation to be submitted by the president to the Board."
$ital4 = $para4.Characters(64, 23)
$ital4.Font.Italic = $true

# --- Edit 5: Slide 9 - "more refined" italic -> bold ---
$s5 = $p.Slides.Item(9)
$tr5 = $s5.Shapes.Item(1).TextFrame.TextRange
$para5 = $tr5.Paragraphs(1,1)
$whole5 = $para5.Characters(1, 138)
$whole5.Text = "3/3/2021 President Fritz writes “proposing a new and more refined governance plan.” Just 20 days were initially given for campus feedback."
$bold5 = $para5.Characters(54, 12)
$bold5.Font.Bold = $true

Write-Host "Edits applied."
